$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving formatting like
# trailing zeros / leading zeros / subscript digits that Excel would
# otherwise "helpfully" reinterpret as a number and reformat.
function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "34.092.69"
Set-TextCell $ws.Range("E2") "  -0.53%  "
Set-TextCell $ws.Range("D3") "1.777.27"
Set-TextCell $ws.Range("E3") "  -2.53%  "
Set-TextCell $ws.Range("E4") "  +0.01%  "
Set-TextCell $ws.Range("D5") "224.87"
Set-TextCell $ws.Range("E5") "  -2.03%  "
Set-TextCell $ws.Range("D6") "0.550"
Set-TextCell $ws.Range("E6") "  +0.56%  "
Set-TextCell $ws.Range("E7") "  +0.01%  "
Set-TextCell $ws.Range("D8") "31.81"
Set-TextCell $ws.Range("E8") "  +0.72%  "
Set-TextCell $ws.Range("E9") "  -1.28%  "
Set-TextCell $ws.Range("E10") "  -2.42%  "
Set-TextCell $ws.Range("E11") "  -0.07%  "
Set-TextCell $ws.Range("D12") "2.032.69"
Set-TextCell $ws.Range("E12") "  -2.56%  "
Set-TextCell $ws.Range("D13") "11.09"
Set-TextCell $ws.Range("E13") "  +6.60%  "
Set-TextCell $ws.Range("D14") "1.776.94"
Set-TextCell $ws.Range("E14") "  -2.55%  "
Set-TextCell $ws.Range("E15") "  -3.38%  "
Set-TextCell $ws.Range("D16") "34.087.94"
Set-TextCell $ws.Range("E16") "  -0.39%  "
Set-TextCell $ws.Range("E17") "  -2.10%  "
Set-TextCell $ws.Range("D18") "68.55"
Set-TextCell $ws.Range("E18") "  -1.77%  "
Set-TextCell $ws.Range("D19") "254.36"
Set-TextCell $ws.Range("E19") "  -1.51%  "
Set-TextCell $ws.Range("D20") "0.0₃0737"
Set-TextCell $ws.Range("E20") "  -2.18%  "
Set-TextCell $ws.Range("D21") "0.999"
Set-TextCell $ws.Range("E21") "  +0.05%  "
Set-TextCell $ws.Range("D22") "10.35"
Set-TextCell $ws.Range("E22") "  -2.32%  "
Set-TextCell $ws.Range("D23") "4.18"
Set-TextCell $ws.Range("E23") "  -3.82%  "
Set-TextCell $ws.Range("E24") "  -3.38%  "
Set-TextCell $ws.Range("D25") "156.97"
Set-TextCell $ws.Range("E25") "  -1.78%  "
Set-TextCell $ws.Range("D26") "16.35"
Set-TextCell $ws.Range("E26") "  -1.63%  "
Set-TextCell $ws.Range("D27") "6.99"
Set-TextCell $ws.Range("E27") "  -2.24%  "
Set-TextCell $ws.Range("E28") "  -1.44%  "
Set-TextCell $ws.Range("E29") "  +0.07%  "
Set-TextCell $ws.Range("D30") "3.76"
Set-TextCell $ws.Range("E30") "  -3.23%  "
Set-TextCell $ws.Range("D31") "0.0512"
Set-TextCell $ws.Range("E31") "  -1.63%  "
Set-TextCell $ws.Range("E32") "  -1.62%  "
Set-TextCell $ws.Range("D33") "3.58"
Set-TextCell $ws.Range("E33") "  +0.51%  "
Set-TextCell $ws.Range("E34") "  +1.96%  "
Set-TextCell $ws.Range("D35") "1.435.85"
Set-TextCell $ws.Range("E35") "  -7.39%  "
Set-TextCell $ws.Range("D36") "1.05"
Set-TextCell $ws.Range("E36") "  -3.86%  "
Set-TextCell $ws.Range("E37") "  -0.93%  "
Set-TextCell $ws.Range("D38") "0.622"
Set-TextCell $ws.Range("E38") "  -1.71%  "
Set-TextCell $ws.Range("E39") "  +1.34%  "
Set-TextCell $ws.Range("D40") "82.68"
Set-TextCell $ws.Range("E40") "  -2.60%  "
Set-TextCell $ws.Range("D41") "2.34"
Set-TextCell $ws.Range("E41") "  +0.55%  "
Set-TextCell $ws.Range("D42") "0.885"
Set-TextCell $ws.Range("E42") "  -3.70%  "
Set-TextCell $ws.Range("D43") "2.05"
Set-TextCell $ws.Range("E43") "  -5.39%  "
Set-TextCell $ws.Range("D44") "0.0513"
Set-TextCell $ws.Range("E44") "  -2.36%  "
Set-TextCell $ws.Range("E45") "  -1.86%  "
Set-TextCell $ws.Range("D46") "1.932.78"
Set-TextCell $ws.Range("E46") "  -2.78%  "
Set-TextCell $ws.Range("E47") "  +0.38%  "
Set-TextCell $ws.Range("D48") "12.20"
Set-TextCell $ws.Range("E48") "  -2.32%  "
Set-TextCell $ws.Range("E49") "  +0.00%  "
Set-TextCell $ws.Range("D50") "98.12"
Set-TextCell $ws.Range("E50") "  +1.08%  "
Set-TextCell $ws.Range("D51") "49.46"
Set-TextCell $ws.Range("E51") "  -6.52%  "
